# Update "想去人数" (F column) counts across the four sheets to match
# the newly generated data snapshot (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1680
$ws.Range("F3").Value  = 9592
$ws.Range("F5").Value  = 788
$ws.Range("F13").Value = 1516
$ws.Range("F15").Value = 318
$ws.Range("F19").Value = 418
$ws.Range("F23").Value = 5
$ws.Range("F29").Value = 613
$ws.Range("F30").Value = 649
$ws.Range("F33").Value = 185
$ws.Range("F34").Value = 91
$ws.Range("F35").Value = 45
$ws.Range("F39").Value = 344
$ws.Range("F40").Value = 637
$ws.Range("F43").Value = 335
$ws.Range("F46").Value = 63

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value  = 73
$ws.Range("F12").Value = 63
$ws.Range("F14").Value = 20
$ws.Range("F23").Value = 666
$ws.Range("F24").Value = 39

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 2431
$ws.Range("F7").Value  = 3779
$ws.Range("F8").Value  = 25
$ws.Range("F10").Value = 143
$ws.Range("F11").Value = 130

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1680
$ws.Range("F4").Value  = 9592
$ws.Range("F7").Value  = 3779
$ws.Range("F8").Value  = 788
$ws.Range("F9").Value  = 143
$ws.Range("F10").Value = 143
$ws.Range("F18").Value = 130
$ws.Range("F19").Value = 1516
$ws.Range("F21").Value = 318
$ws.Range("F33").Value = 613
$ws.Range("F34").Value = 649
$ws.Range("F35").Value = 39
$ws.Range("F38").Value = 45
$ws.Range("F41").Value = 344
$ws.Range("F43").Value = 637
$ws.Range("F46").Value = 335
